$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 9).Value = 'b'
$ws.Cells.Item(2, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(6, 9).Value = 'sv'
$ws.Cells.Item(6, 10).Value = 'Statement-opinion'
$ws.Cells.Item(16, 9).Value = '%'
$ws.Cells.Item(16, 10).Value = 'Uninterpretable'
$ws.Cells.Item(30, 9).Value = 'sd'
$ws.Cells.Item(30, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(44, 9).Value = 'b'
$ws.Cells.Item(44, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(51, 9).Value = 'sd'
$ws.Cells.Item(51, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(53, 9).Value = '%'
$ws.Cells.Item(53, 10).Value = 'Uninterpretable'
$ws.Cells.Item(61, 9).Value = 'sd'
$ws.Cells.Item(61, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(63, 9).Value = '%'
$ws.Cells.Item(63, 10).Value = 'Uninterpretable'
$ws.Cells.Item(74, 9).Value = 'qy'
$ws.Cells.Item(74, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(86, 9).Value = 'sd'
$ws.Cells.Item(86, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(90, 9).Value = 'sv'
$ws.Cells.Item(90, 10).Value = 'Statement-opinion'
$ws.Cells.Item(94, 9).Value = 'sv'
$ws.Cells.Item(94, 10).Value = 'Statement-opinion'
$ws.Cells.Item(97, 9).Value = 'b'
$ws.Cells.Item(97, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(102, 9).Value = 'ba'
$ws.Cells.Item(102, 10).Value = 'Appreciation'
$ws.Cells.Item(115, 9).Value = 'aa'
$ws.Cells.Item(115, 10).Value = 'Agree/Accept'
$ws.Cells.Item(116, 9).Value = 'sv'
$ws.Cells.Item(116, 10).Value = 'Statement-opinion'
$ws.Cells.Item(117, 9).Value = 'sd'
$ws.Cells.Item(117, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(126, 9).Value = 'sv'
$ws.Cells.Item(126, 10).Value = 'Statement-opinion'
$ws.Cells.Item(128, 9).Value = 'sv'
$ws.Cells.Item(128, 10).Value = 'Statement-opinion'
$ws.Cells.Item(129, 9).Value = 'ba'
$ws.Cells.Item(129, 10).Value = 'Appreciation'
$ws.Cells.Item(139, 9).Value = 'sv'
$ws.Cells.Item(139, 10).Value = 'Statement-opinion'
$ws.Cells.Item(151, 9).Value = 'b'
$ws.Cells.Item(151, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(157, 9).Value = '%'
$ws.Cells.Item(157, 10).Value = 'Uninterpretable'
$ws.Cells.Item(161, 9).Value = 'sv'
$ws.Cells.Item(161, 10).Value = 'Statement-opinion'
$ws.Cells.Item(172, 9).Value = 'sv'
$ws.Cells.Item(172, 10).Value = 'Statement-opinion'
$ws.Cells.Item(173, 9).Value = 'sv'
$ws.Cells.Item(173, 10).Value = 'Statement-opinion'
$ws.Cells.Item(178, 9).Value = 'aa'
$ws.Cells.Item(178, 10).Value = 'Agree/Accept'
$ws.Cells.Item(198, 9).Value = 'sv'
$ws.Cells.Item(198, 10).Value = 'Statement-opinion'
$ws.Cells.Item(199, 9).Value = 'sv'
$ws.Cells.Item(199, 10).Value = 'Statement-opinion'
$ws.Cells.Item(206, 9).Value = 'aa'
$ws.Cells.Item(206, 10).Value = 'Agree/Accept'
$ws.Cells.Item(210, 9).Value = 'ba'
$ws.Cells.Item(210, 10).Value = 'Appreciation'
$ws.Cells.Item(233, 9).Value = 'aa'
$ws.Cells.Item(233, 10).Value = 'Agree/Accept'
$ws.Cells.Item(238, 9).Value = 'sv'
$ws.Cells.Item(238, 10).Value = 'Statement-opinion'
$ws.Cells.Item(239, 9).Value = 'sd'
$ws.Cells.Item(239, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(241, 9).Value = 'sd'
$ws.Cells.Item(241, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(270, 9).Value = 'b'
$ws.Cells.Item(270, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(274, 9).Value = 'sd'
$ws.Cells.Item(274, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(302, 9).Value = 'sd'
$ws.Cells.Item(302, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(303, 9).Value = 'ba'
$ws.Cells.Item(303, 10).Value = 'Appreciation'
$ws.Cells.Item(305, 9).Value = 'ba'
$ws.Cells.Item(305, 10).Value = 'Appreciation'
$ws.Cells.Item(333, 9).Value = 'sd'
$ws.Cells.Item(333, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(334, 9).Value = 'sv'
$ws.Cells.Item(334, 10).Value = 'Statement-opinion'
$ws.Cells.Item(348, 9).Value = 'sd'
$ws.Cells.Item(348, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(351, 9).Value = 'sv'
$ws.Cells.Item(351, 10).Value = 'Statement-opinion'
$ws.Cells.Item(368, 9).Value = 'sd'
$ws.Cells.Item(368, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(372, 9).Value = 'sv'
$ws.Cells.Item(372, 10).Value = 'Statement-opinion'
$ws.Cells.Item(373, 9).Value = 'sd'
$ws.Cells.Item(373, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(376, 9).Value = 'sv'
$ws.Cells.Item(376, 10).Value = 'Statement-opinion'
$ws.Cells.Item(412, 9).Value = 'sd'
$ws.Cells.Item(412, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(414, 9).Value = 'sv'
$ws.Cells.Item(414, 10).Value = 'Statement-opinion'
$ws.Cells.Item(446, 9).Value = 'sd'
$ws.Cells.Item(446, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(461, 9).Value = 'sv'
$ws.Cells.Item(461, 10).Value = 'Statement-opinion'
$ws.Cells.Item(475, 9).Value = 'sv'
$ws.Cells.Item(475, 10).Value = 'Statement-opinion'
$ws.Cells.Item(484, 9).Value = 'sv'
$ws.Cells.Item(484, 10).Value = 'Statement-opinion'
$ws.Cells.Item(487, 9).Value = 'sd'
$ws.Cells.Item(487, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(490, 9).Value = 'sd'
$ws.Cells.Item(490, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(494, 9).Value = 'aa'
$ws.Cells.Item(494, 10).Value = 'Agree/Accept'
$ws.Cells.Item(499, 9).Value = 'b'
$ws.Cells.Item(499, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(503, 9).Value = 'aa'
$ws.Cells.Item(503, 10).Value = 'Agree/Accept'
$ws.Cells.Item(511, 9).Value = 'sd'
$ws.Cells.Item(511, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(519, 9).Value = 'aa'
$ws.Cells.Item(519, 10).Value = 'Agree/Accept'
$ws.Cells.Item(520, 9).Value = '%'
$ws.Cells.Item(520, 10).Value = 'Uninterpretable'
$ws.Cells.Item(521, 9).Value = 'sd'
$ws.Cells.Item(521, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(522, 9).Value = 'sv'
$ws.Cells.Item(522, 10).Value = 'Statement-opinion'
$ws.Cells.Item(524, 9).Value = 'sv'
$ws.Cells.Item(524, 10).Value = 'Statement-opinion'
$ws.Cells.Item(549, 9).Value = 'sd'
$ws.Cells.Item(549, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(550, 9).Value = 'sd'
$ws.Cells.Item(550, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(556, 9).Value = 'sv'
$ws.Cells.Item(556, 10).Value = 'Statement-opinion'
$ws.Cells.Item(557, 9).Value = 'sv'
$ws.Cells.Item(557, 10).Value = 'Statement-opinion'
$ws.Cells.Item(559, 9).Value = 'qy'
$ws.Cells.Item(559, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(560, 9).Value = 'sd'
$ws.Cells.Item(560, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(577, 9).Value = 'sd'
$ws.Cells.Item(577, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(583, 9).Value = 'sv'
$ws.Cells.Item(583, 10).Value = 'Statement-opinion'
$ws.Cells.Item(586, 9).Value = 'sd'
$ws.Cells.Item(586, 10).Value = 'Statement-non-opinion'
